$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.957.82'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.265.05'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.12'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.87'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.08'
$ws.Range('E10').Value = '  +7.29%  '
$ws.Range('E11').Value = '  -1.11%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.64'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '2.614.43'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.40'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '2.265.93'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('D18').Value = '41.862.35'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('E19').Value = '  -3.17%  '
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.88'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.30'
$ws.Range('E23').Value = '  -2.79%  '
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.93'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.66'
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.71'
$ws.Range('E28').Value = '  +5.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.50'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.40'
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0737'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.104'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('E40').Value = '  -2.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.02'
$ws.Range('E41').Value = '  +2.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.38'
$ws.Range('E42').Value = '  +5.75%  '
$ws.Range('D43').Value = '1.977.59'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.74'
$ws.Range('E45').Value = '  -4.76%  '
$ws.Range('E46').Value = '  +0.87%  '
$ws.Range('E47').Value = '  -3.86%  '
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '72.97'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.80'
$ws.Range('E51').Value = '  -1.25%  '
